$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 21.31228666666667
$ws.Range("H2").Value = 63.93686
$ws.Range("I2").Value = 0.9506775731819035
$ws.Range("J2").Value = 0.9506775731819034
$ws.Range("M2").Value = 0.029424
$ws.Range("N2").Value = 0.08827199999999999
$ws.Range("O2").Value = 0.1473063425232919
$ws.Range("P2").Value = 0.1473063425232919
$ws.Range("Q2").Value = 0.62709272288
$ws.Range("R2").Value = 5.643834505919999
$ws.Range("S2").Value = 0.1400408362243454
$ws.Range("T2").Value = 0.1400408362243454

# Row 3
$ws.Range("G3").Value = 21.31228666666667
$ws.Range("H3").Value = 63.93686
$ws.Range("I3").Value = 0.9506775731819035
$ws.Range("J3").Value = 0.9506775731819034
$ws.Range("O3").Value = 0.852693657476708
$ws.Range("P3").Value = 0.852693657476708
$ws.Range("Q3").Value = 3.629972601926667
$ws.Range("R3").Value = 32.66975341734
$ws.Range("S3").Value = 0.810636736957558
$ws.Range("T3").Value = 0.8106367369575579

# Row 4
$ws.Range("G4").Value = 0.7500946666666666
$ws.Range("I4").Value = 0.03345948693899053
$ws.Range("J4").Value = 0.03345948693899053
$ws.Range("M4").Value = 0.029424
$ws.Range("N4").Value = 0.08827199999999999
$ws.Range("O4").Value = 0.1473063425232919
$ws.Range("P4").Value = 0.1473063425232919
$ws.Range("Q4").Value = 0.02207078547199999
$ws.Range("R4").Value = 0.198637069248
$ws.Range("S4").Value = 0.004928794643688551
$ws.Range("T4").Value = 0.004928794643688551

# Row 5
$ws.Range("G5").Value = 0.7500946666666666
$ws.Range("I5").Value = 0.03345948693899053
$ws.Range("J5").Value = 0.03345948693899053
$ws.Range("O5").Value = 0.852693657476708
$ws.Range("P5").Value = 0.852693657476708
$ws.Range("S5").Value = 0.02853069229530197
$ws.Range("T5").Value = 0.02853069229530197

# Row 6
$ws.Range("I6").Value = 0.01586293987910606
$ws.Range("J6").Value = 0.01586293987910605
$ws.Range("M6").Value = 0.029424
$ws.Range("N6").Value = 0.08827199999999999
$ws.Range("O6").Value = 0.1473063425232919
$ws.Range("P6").Value = 0.1473063425232919
$ws.Range("Q6").Value = 0.010463625568
$ws.Range("R6").Value = 0.09417263011199999
$ws.Range("S6").Value = 0.002336711655257984
$ws.Range("T6").Value = 0.002336711655257984

# Row 7
$ws.Range("I7").Value = 0.01586293987910606
$ws.Range("J7").Value = 0.01586293987910605
$ws.Range("O7").Value = 0.852693657476708
$ws.Range("P7").Value = 0.852693657476708
$ws.Range("S7").Value = 0.01352622822384807
$ws.Range("T7").Value = 0.01352622822384807
